$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = 4

# Row 2 (existing row, values updated - shared-string indices shift but text stays same content)
$ws.Range("A2").Value = "MAUPASSANT"
$ws.Range("B2").Value = "Guy"
$ws.Range("C2").Value = "m"
$ws.Range("D2").Value = "Le Horla"
$ws.Range("E2").Value = "Lois"
$ws.Range("G2").Value = 625859696
$ws.Range("H2").Value = "Les Lillas"
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = "NaN"
$ws.Range("K2").Value = "F"
$ws.Range("L2").Value = "O"
$ws.Range("N2").Value = 2500
$ws.Range("O2").Value = "FR78"

# Row 3 (new data row)
$ws.Range("A3").Value = "A"
$ws.Range("B3").Value = "b"
$ws.Range("C3").Value = "f"
$ws.Range("D3").Value = "a"
$ws.Range("E3").Value = "a"
$ws.Range("F3").Value = 36145
$ws.Range("G3").Value = 654987852
$ws.Range("H3").Value = "Les Lillas"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "F"
$ws.Range("L3").Value = "O"
$ws.Range("M3").Value = 36145
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = "fr"

# Row 4 (new data row)
$ws.Range("A4").Value = "SDF"
$ws.Range("B4").Value = "f"
$ws.Range("C4").Value = "f"
$ws.Range("D4").Value = "qdsf"
$ws.Range("E4").Value = "qsdf"
$ws.Range("F4").Value = 36145
$ws.Range("G4").Value = 2457865424
$ws.Range("H4").Value = "Les Lillas"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = "F"
$ws.Range("L4").Value = "O"
$ws.Range("M4").Value = 36145
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = "fr"

# Row 5 (new data row)
$ws.Range("A5").Value = "AZE"
$ws.Range("B5").Value = "aze"
$ws.Range("C5").Value = "f"
$ws.Range("D5").Value = "aze"
$ws.Range("E5").Value = "e"
$ws.Range("F5").Value = 36145
$ws.Range("G5").Value = 245789654
$ws.Range("H5").Value = "Les Lillas"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = "F"
$ws.Range("L5").Value = "F"

# Copy date formatting (style index already used by F2/M2/H2) down onto
# the new date cells without introducing new number-format styles.
$ws.Range("F2").Copy()
$ws.Range("F3:F5").PasteSpecial(-4122)
$ws.Range("M2").Copy()
$ws.Range("M3:M4").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Clear the lingering cell selection shown in the sheet view
$ws.Range("A1").Select()
